# Add a leading "id" column (1..61) to the groceries table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a blank column at the far left; this pushes the existing
#    "name, price, contains_*" data (and their column formatting) one
#    column to the right (A->B, B->C, ... H->I).
$ws.Columns.Item(1).Insert()

# 2. New header + sequential numeric ids for each data row.
$ws.Range("A1").Value = "id"
for ($r = 2; $r -le 62; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# 3. Give the new id column a narrow, auto-fit-like width. (Column B keeps
#    the width it inherited from the original column A untouched.)
$ws.Columns.Item(1).ColumnWidth = 3.6666666666666665

# 4. The table needs to grow by one column on its left edge. This engine's
#    ListObject tracking does not follow a plain column insert, so rebuild
#    the table over the new full range to get a clean column/header bind.
$tbl = $ws.ListObjects.Item(1)
$tbl.Unlist()
$newTbl = $ws.ListObjects.Add(1, $ws.Range("A1:I62"), $null, 1)
$newTbl.Name = "Table1"

# 5. Restore the selection shown in the saved workbook.
$ws.Range("C5").Select()
